# Fruta / hortaliza, semanal
# Insert two new price-report rows (895, 896) for "Lane Late" oranges
# (a new weekly reporting date, serial 44939) right before the existing
# "Navel Late" block that used to start at row 895. Everything that was
# at rows 895..936 shifts down to 897..938.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 895, pushing the old 895..936 block down to 897..938.
$ws.Rows.Item(895).Insert()
$ws.Rows.Item(895).Insert()

# --- New row 895: Lane Late / Primera, fecha 44939 ---
$ws.Cells.Item(895, 1).Value = 8
$ws.Cells.Item(895, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(895, 3).Value = "Coquimbo"
$ws.Cells.Item(895, 4).Value = 44939
$ws.Cells.Item(895, 5).Value = 4
$ws.Cells.Item(895, 6).Value = "Fruta"
$ws.Cells.Item(895, 7).Value = 100102
$ws.Cells.Item(895, 8).Value = "Cítricos"
$ws.Cells.Item(895, 9).Value = 100102005
$ws.Cells.Item(895, 10).Value = "Naranja"
$ws.Cells.Item(895, 11).Value = "Lane Late"
$ws.Cells.Item(895, 12).Value = "Primera"
$ws.Cells.Item(895, 13).Value = 16
$ws.Cells.Item(895, 14).Value = 250000
$ws.Cells.Item(895, 15).Value = 260000
$ws.Cells.Item(895, 16).Value = 255000
$ws.Cells.Item(895, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(895, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(895, 19).Value = 638
$ws.Cells.Item(895, 20).Value = 400

# --- New row 896: Lane Late / Segunda, fecha 44939 ---
$ws.Cells.Item(896, 1).Value = 8
$ws.Cells.Item(896, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(896, 3).Value = "Coquimbo"
$ws.Cells.Item(896, 4).Value = 44939
$ws.Cells.Item(896, 5).Value = 4
$ws.Cells.Item(896, 6).Value = "Fruta"
$ws.Cells.Item(896, 7).Value = 100102
$ws.Cells.Item(896, 8).Value = "Cítricos"
$ws.Cells.Item(896, 9).Value = 100102005
$ws.Cells.Item(896, 10).Value = "Naranja"
$ws.Cells.Item(896, 11).Value = "Lane Late"
$ws.Cells.Item(896, 12).Value = "Segunda"
$ws.Cells.Item(896, 13).Value = 10
$ws.Cells.Item(896, 14).Value = 210000
$ws.Cells.Item(896, 15).Value = 220000
$ws.Cells.Item(896, 16).Value = 215000
$ws.Cells.Item(896, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(896, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(896, 19).Value = 538
$ws.Cells.Item(896, 20).Value = 400
